$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) for every data row (2-31) was stored as the
# malformed text "4-29-2011-12" because of how the NBA stats source
# displayed the date. Correct it to the proper text "2012-04-29".
$oldValue = "4-29-2011-12"
$newValue = "2012-04-29"

$rng = $ws.Range("BF2:BF31")

# Switch the range to Text format first so that assigning the
# date-looking string "2012-04-29" is stored as literal text instead of
# being auto-converted by Excel into a date serial number.
$rng.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}

# Restore the default "Normal" style on the range so the cells keep the
# same look as before (no lingering custom text number format).
$rng.Style = "Normal"
